$wb = $excel.ActiveWorkbook
$wsLogs = $wb.Worksheets.Item("Logs")
$wsDash = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: append row 17 -------------------------------------------------
$wsLogs.Range("A17").Value = "Weten jullie al iets over mijn retour?"
$wsLogs.Range("B17").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C17").Value = "Testmail #1: Weten jullie al iets over mijn retour?"
$wsLogs.Range("D17").Value = "Retour / Terugbetaling"
$wsLogs.Range("E17").Value = "Beste klant,`r`nBedankt voor je bericht. Om je retour te kunnen bekijken, hebben we wat meer informatie nodig. Zou je ons je bestelnummer kunnen doorgeven, zodat we de status van je retour kunnen controleren?`r`nAlvast bedankt voor je medewerking.`r`nMet vriendelijke groet,`r`n[Naam bedrijf]"
$wsLogs.Range("F17").Value = "2025-08-06 20:33:19"
$wsLogs.Range("G17").Value = "Ja"
$wsLogs.Range("H17").Value = "Nee"
$wsLogs.Range("I17").Value = "Ja"
$wsLogs.Range("J17").Value = "Nee"
$wsLogs.Rows(17).AutoFit()

# --- Logs sheet: expand conditional formatting ranges to include row 17 -------
function Expand-ConditionalFormatting($column) {
    $oldRange = $wsLogs.Range($column + "2:" + $column + "16")
    $newRange = $wsLogs.Range($column + "2:" + $column + "17")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

Expand-ConditionalFormatting "D"
Expand-ConditionalFormatting "G"
Expand-ConditionalFormatting "H"
Expand-ConditionalFormatting "I"
Expand-ConditionalFormatting "J"

# --- Dashboard sheet: append row 6 ---------------------------------------------
$wsDash.Range("A6").Value = "Retour / Terugbetaling"
$wsDash.Range("B6").Value = 1

# --- Dashboard chart: extend series category/value references to row 6 --------
$chart = $wsDash.ChartObjects(1).Chart
$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$6,'Dashboard'!`$B`$2:`$B`$6,1)"
